$wb = $excel.ActiveWorkbook

# Sheet "5" (sheet1.xml) holds the main trades table (rows 2-9).
# The first trade row (row 2, ID=7 / USDBTC / buy) was removed and the
# rows below it shifted up by one.
$ws5 = $wb.Worksheets.Item("5")
$ws5.Rows.Item(2).Delete()

# Sheet "8" (sheet4.xml) has a single trade row whose COMMENTS cell (K2)
# text changed from "strategy 2 comment" to "strategy 3 comment".
$ws8 = $wb.Worksheets.Item("8")
$ws8.Range("K2").Value = "strategy 3 comment"
